# Apply the author's edits to the "planets" sheet (Sheet1):
#  - Undare's (row 14) Moons count changed from 7 to 73
#  - Active selection left on E14 (where the edit was made)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the data value: Moons for Undare (row 14, column B) 7 -> 73
$ws.Range("B14").Value = 73

# Leave the selection where the edit was made (E14), matching the saved view state
$ws.Range("E14").Select()
